# Apply the "Holden scheme" update to the UniformF-HW03 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relabel rows 16-19 (HKL index 14-17) from the HexGrid-90degTilt* labels
#    to the new Holden* labels.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 2. Reorder the HKL-triplet header labels in row 2, columns C:J.
$ws.Range("C2").Value = "[3, 2, 1]"
$ws.Range("D2").Value = "[3, 1, 0]"
$ws.Range("E2").Value = "[2, 2, 2]"
$ws.Range("F2").Value = "[1, 1, 0]"
$ws.Range("G2").Value = "[2, 0, 0]"
$ws.Range("H2").Value = "[2, 2, 0]"
$ws.Range("I2").Value = "[4, 0, 0]"
$ws.Range("J2").Value = "[2, 1, 1]"

# 3. Append four new rows (20-23) holding the HexGrid data that used to
#    live in rows 16-19, each filled with 1s across C:T, matching the
#    existing row layout/style.
$holdenHexRows = @(
    @{ Row = 20; Idx = 18; Label = "HexGrid-90degTilt2.5degRes" },
    @{ Row = 21; Idx = 19; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 22; Idx = 20; Label = "HexGrid-90degTilt10degRes" },
    @{ Row = 23; Idx = 21; Label = "HexGrid-90degTilt15degRes" }
)

foreach ($entry in $holdenHexRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Idx
    $ws.Range("A" + ($r - 1)).Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Cells.Item($r, 2).Value = $entry.Label
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
